{"js": "// Word homework document edit:\n// 1. Highlight \"bin has a different number\" green.\n// 2. Highlight \"parking spot\" cyan.\n// 3. Move the \"_GoBack\" bookmark from just before \" is also recorded\"\n//    (after \"bin this date\") to just before \"parking spot\"\n//    (i.e. right after \"Some bays have a \").\n\nconst body = context.document.body;\n\n// --- 1. Highlight \"bin has a different number\" in green ---\nconst greenResults = body.search(\"bin has a different number\", { matchCase: true });\ngreenResults.load(\"text\");\nawait context.sync();\nif (greenResults.items.length > 0) {\n  greenResults.items[0].font.highlightColor = \"green\";\n}\n\n// --- 2. Highlight \"parking spot\" in cyan ---\nconst cyanResults = body.search(\"parking spot\", { matchCase: true });\ncyanResults.load(\"text\");\nawait context.sync();\nif (cyanResults.items.length > 0) {\n  cyanResults.items[0].font.highlightColor = \"cyan\";\n}\nawait context.sync();\n\n// --- 3. Relocate the \"_GoBack\" bookmark ---\n// Remove it from its original position (it will be re-inserted below).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-search, since the document changed, to get a fresh anchor for \"parking spot\".\nconst parkingResults = body.search(\"parking spot\", { matchCase: true });\nparkingResults.load(\"text\");\nawait context.sync();\nif (parkingResults.items.length > 0) {\n  const startRange = parkingResults.items[0].getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word homework document edit (COM/PowerShell):\n# 1. Highlight \"bin has a different number\" green.\n# 2. Highlight \"parking spot\" cyan.\n# 3. Move the \"_GoBack\" bookmark from just before \" is also recorded\"\n#    (after \"bin this date\") to just before \"parking spot\"\n#    (i.e. right after \"Some bays have a \").\n\n$d = $word.ActiveDocument\n\n# --- 1. Highlight \"bin has a different number\" in green ---\n$rngGreen = $d.Content\n$rngGreen.Find.Execute(\"bin has a different number\") | Out-Null\n$rngGreen.Font.HighlightColorIndex = \"wdBrightGreen\"\n\n# --- 2. Highlight \"parking spot\" in cyan ---\n$rngCyan = $d.Content\n$rngCyan.Find.Execute(\"parking spot\") | Out-Null\n$rngCyan.Font.HighlightColorIndex = \"wdTurquoise\"\n\n# --- 3. Relocate the \"_GoBack\" bookmark ---\n# Remove it from its original position (it will be re-inserted below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-find \"parking spot\" to get a fresh range/start position, then drop a\n# zero-length bookmark right before it (i.e. right after \"Some bays have a \").\n$rngBookmark = $d.Content\n$rngBookmark.Find.Execute(\"parking spot\") | Out-Null\n$startRng = $d.Range($rngBookmark.Start, $rngBookmark.Start)\n$d.Bookmarks.Add(\"_GoBack\", $startRng)\n"}
